$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column C (LFP_decay) equal to column B (general decay) for rows 2..102
for ($r = 2; $r -le 102; $r++) {
    $ws.Cells.Item($r, 3).Value2 = $ws.Cells.Item($r, 2).Value2
}

# Update the active selection to D2 (as recorded in the saved file)
$ws.Range("D2").Select()
